$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows right after the header row (row 1). This pushes the
# existing 2021-2024 data rows down from rows 2-5 to rows 5-8.
$ws.Rows("2:4").Insert()

# Row insertion copies the formatting of the row above (the header row),
# which is not what the data rows should look like - strip it back off the
# new B:D cells so they stay unstyled like the other data rows.
$ws.Range("B2:D4").ClearFormats()

# Give the new index cells (column A) the same style as the existing index
# column cells.
$ws.Range("A5").Copy()
$ws.Range("A2:A4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# New data for 2018, 2019, 2020 in the freshly inserted rows 2-4.
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "2018"
$ws.Range("C2").Value = 1871
$ws.Range("D2").Value = 32.23

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "2019"
$ws.Range("C3").Value = 1260
$ws.Range("D3").Value = 21.71

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "2020"
$ws.Range("C4").Value = 682
$ws.Range("D4").Value = 11.75

# Re-number the index column (A) sequentially 0..6 across all data rows,
# since the pre-existing rows (now 5-8) shifted down.
$ws.Range("A5").Value = 3
$ws.Range("A6").Value = 4
$ws.Range("A7").Value = 5
$ws.Range("A8").Value = 6
